$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 265; this pushes the existing row 265 (and everything
# below it, through the old row 382) down by one row, to row 266..383.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new record.
$ws.Range("A265").Value = 4
$ws.Range("B265").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C265").Value = "Los Lagos"
$ws.Range("D265").Value = 44917
$ws.Range("D265").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E265").Value = 10
$ws.Range("F265").Value = 100112003
$ws.Range("G265").Value = "Ajo"
$ws.Range("H265").Value = "Chino"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 70
$ws.Range("K265").Value = 20000
$ws.Range("L265").Value = 20000
$ws.Range("M265").Value = 20000
$ws.Range("N265").Value = "$/caja 10 kilos"
$ws.Range("O265").Value = "China"
$ws.Range("P265").Value = 2000
$ws.Range("Q265").Value = 10
$ws.Range("R265").Value = "Hortaliza"
